$wb = $excel.ActiveWorkbook

# --- "Insumos" sheet: clear the leftover test-name row (luis / estela / jimenez) ---
$wsInsumos = $wb.Worksheets.Item("Insumos")
$wsInsumos.Range("A7:C7").ClearContents()

# --- "Procesados" sheet: row heights settle to 12.75 and the view selection moves ---
$wsProcesados = $wb.Worksheets.Item("Procesados")
$wsProcesados.Rows("1:3").RowHeight = 12.75
$wsProcesados.Activate()
$wsProcesados.Range("F13:F14").Select()

# --- "Insumos" becomes the active sheet/tab (it was "NoProcesados" before) ---
$wsInsumos.Activate()
$wsInsumos.Range("F14").Select()
